$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the entry date value from 17.12.2023 to 01.01.2025 (text value in B6)
$ws.Range("B6").Value = "01.01.2025"

# Move active selection to B7 (mirrors author's last-saved cursor position)
$ws.Range("B7").Select()
